# Updates cryptos list: refresh Price (D) / Volume(1h) (E) figures, and
# correct the B/C/D/E data for rows 26-28 (Binance-PegBSC-USD, Kaspa, PEPE).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.271.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "2.991.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'506.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").Value = "'137.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("D12").Value = "3.504.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").Value = "56.248.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").Value = "2.991.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D19").Value = "'12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "

$ws.Range("D20").Value = "'8.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "

$ws.Range("D21").Value = "'331.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'0.495"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").Value = "'64.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.76%  "

$ws.Range("D25").Value = "3.117.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0940"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.28%  "

$ws.Range("E29").Value = "  -3.49%  "

$ws.Range("D30").Value = "'6.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "'20.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").Value = "'152.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.69%  "

$ws.Range("D35").Value = "'4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.31%  "

$ws.Range("D36").Value = "'5.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("D37").Value = "'26.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.70%  "

$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").Value = "'0.0661"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "3.029.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "'36.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  +1.03%  "

$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("D45").Value = "2.188.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").Value = "'5.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("D48").Value = "'0.923"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("D49").Value = "'0.0236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'19.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").Value = "'0.0852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
